$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for a handful of rows per the repull/mean
# recalculation described in the commit message.
$ws.Range("F3").Value = 4
$ws.Range("F4").Value = -3
$ws.Range("F7").Value = -8
$ws.Range("F9").Value = -5
